# Re-pick two variables for the example2 dataset ("Data" sheet) + refresh
# the matching rows in the "Codebook" sheet. The old "Race"/"BMI" columns
# are replaced with "Preferred exercise type" / "Resting Heart Rate".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Codebook")

# ---------------------------------------------------------------------
# "Data" sheet: column D (exercise type) and column E (resting heart rate)
# ---------------------------------------------------------------------
$ws1.Range("D1").Value = "Preferred exercise type"
$ws1.Range("E1").Value = "Resting Heart Rate"

$exercise = @{
    2  = "Strength training"
    3  = "Cardio"
    4  = "None"
    5  = "Strength training"
    6  = "Strength training"
    7  = "Yoga"
    8  = "None"
    9  = "None"
    10 = "Cardio"
    11 = "None"
    12 = "Yoga"
    13 = "Yoga"
    14 = "Cardio"
    15 = "Yoga"
}
$heartrate = @{
    2  = 60
    3  = 65
    4  = 153
    5  = 72
    6  = 85
    7  = 40
    8  = 70
    9  = 90
    10 = 56
    11 = 67
    12 = 60
    13 = 64
    14 = 78
    15 = 60
}

foreach ($r in 2..15) {
    $ws1.Range("D$r").Value = $exercise[$r]
    $ws1.Range("E$r").Value = $heartrate[$r]
}

# column widths for the rewritten columns (closest attainable to the
# authored 19.5 / 22.6640625 character widths)
$ws1.Columns.Item(4).ColumnWidth = 18.585
$ws1.Columns.Item(5).ColumnWidth = 21.75

# ---------------------------------------------------------------------
# "Codebook" sheet: rows 5 & 6 describe the two swapped-in variables
# ---------------------------------------------------------------------
$ws2.Range("A5").Value = "Preferred exercise Type"
$ws2.Range("B5").Value = "Preferred exercise type (Cardio, Strenght training, Yoga, None)"
$ws2.Range("C5").Value = "Cardio/Strength training /Yoga/None"

$ws2.Range("A6").Value = "Resting Heart Rate"
$ws2.Range("B6").Value = "Heart rate when resting ranging from 40 to 153"
$ws2.Range("C6").Value = "numeric value >0 "

$ws2.Columns.Item(1).ColumnWidth = 19.92
$ws2.Columns.Item(2).ColumnWidth = 47.085
$ws2.Columns.Item(3).ColumnWidth = 42.25

# ---------------------------------------------------------------------
# view state: zoom + selection on both sheets, keep "Data" the active tab
# ---------------------------------------------------------------------
[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 173
$ws2.Range("B5").Select() | Out-Null

[void]$ws1.Activate()
$ws1.Range("D1").Select() | Out-Null

"done"
